$wb = $excel.ActiveWorkbook

# Add "Vertical" sheet right after the existing (last) sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVertical = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsVertical.Name = "Vertical"
$wsVertical.Range("B1").Value = "Vertical"
$wsVertical.Range("A3").Value = "Standard Deviation"
$wsVertical.Range("B3").Value = 0.04829447106908376
$wsVertical.Range("A4").Value = "Maximum"
$wsVertical.Range("B4").Value = 0.3646477683989286

# Add "Updated Lin" sheet right after "Vertical".
$wsUpdatedLin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsVertical)
$wsUpdatedLin.Name = "Updated Lin"
$wsUpdatedLin.Range("B1").Value = "Updated Lin (Z)"
$wsUpdatedLin.Range("A3").Value = "Standard Deviation"
$wsUpdatedLin.Range("B3").Value = 0.04829447106908376
$wsUpdatedLin.Range("A4").Value = "Maximum"
$wsUpdatedLin.Range("B4").Value = 0.3646477683989286
